$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match (bold, border, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-20
$data = @(
    @(2, 5),
    @(6, 9),
    @(2, 5),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 3),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(6, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
